# Natmi LR-pair output: "following Dr Hou advice"
#
# The sending-cluster "ECs" / "FAPs" -> Tnf/Tnfrsf21 -> target-cluster rows are
# recomputed (new average/total expression + specificity numbers), and the
# analysis is extended to also report the reverse direction, where "FAPs" is
# the sending cluster. This doubles the 3 original data rows (ECs -> ECs /
# FAPs / sCs) into 6 rows: the original 3 (now with updated numbers) plus 3
# new ones for FAPs -> ECs / FAPs / sCs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One row per (Sending cluster, Ligand, Receptor, Target cluster) combination,
# followed by the 16 numeric metric columns E..T.
$rows = @(
    @("ECs",  "Tnf", "Tnfrsf21", "ECs",  3, 1, 459.5553626666667, 1378.666088,       0.997342491800172,   0.997342491800172,   3, 1, 12.98318866666667, 38.949566,  0.2170149059254416,  0.2170149059254416,  5966.493976279758, 53698.44578651781, 0.2164381870334598,   0.2164381870334598),
    @("ECs",  "Tnf", "Tnfrsf21", "FAPs", 3, 1, 459.5553626666667, 1378.666088,       0.997342491800172,   0.997342491800172,   3, 1, 2.096730333333333, 6.290191,   0.03504699405682875, 0.03504699405682875, 963.5636687492009, 8672.073018742807, 0.03495385638274341,  0.03495385638274341),
    @("ECs",  "Tnf", "Tnfrsf21", "sCs",  3, 1, 459.5553626666667, 1378.666088,       0.997342491800172,   0.997342491800172,   3, 1, 44.746334,         134.239002, 0.7479381000177296,  0.7479381000177296,  20563.4177493738,  185070.7597443642, 0.7459504483839687,   0.7459504483839687),
    @("FAPs", "Tnf", "Tnfrsf21", "ECs",  2, 0.6666666666666666, 1.224526333333333, 3.673579, 0.002657508199827995, 0.002657508199827995, 3, 1, 12.98318866666667, 38.949566,  0.2170149059254416,  0.2170149059254416,  15.89825641296822, 143.084307716714,  0.000576718891981762, 0.000576718891981762),
    @("FAPs", "Tnf", "Tnfrsf21", "FAPs", 2, 0.6666666666666666, 1.224526333333333, 3.673579, 0.002657508199827995, 0.002657508199827995, 3, 1, 2.096730333333333, 6.290191,   0.03504699405682875, 0.03504699405682875, 2.567501507065445, 23.107513563589,   0.00009313767408534543, 0.00009313767408534543),
    @("FAPs", "Tnf", "Tnfrsf21", "sCs",  2, 0.6666666666666666, 1.224526333333333, 3.673579, 0.002657508199827995, 0.002657508199827995, 3, 1, 44.746334,         134.239002, 0.7479381000177296,  0.7479381000177296,  54.79306430312867, 493.137578728158,  0.001987651633760888, 0.001987651633760888)
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $ws.Range("$($columns[$j])$r").Value = $data[$j]
    }
}
